$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.695.48'
$ws.Range("E2").Value = '  -0.78%  '
$ws.Range("D3").Value = '1.590.87'
$ws.Range("E3").Value = '  -2.41%  '
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").Value = "'" + '208.38'
$ws.Range("E5").Value = '  -1.65%  '
$ws.Range("E6").Value = '  -2.98%  '
$ws.Range("E7").Value = '  +0.17%  '
$ws.Range("D8").Value = "'" + '22.29'
$ws.Range("E8").Value = '  -4.08%  '
$ws.Range("E9").Value = '  -2.02%  '
$ws.Range("D10").Value = "'" + '0.0593'
$ws.Range("E10").Value = '  -2.45%  '
$ws.Range("D11").Value = "'" + '0.0867'
$ws.Range("E11").Value = '  -1.76%  '
$ws.Range("D12").Value = '1.816.58'
$ws.Range("E12").Value = '  -2.40%  '
$ws.Range("D13").Value = '1.590.33'
$ws.Range("E13").Value = '  -2.46%  '
$ws.Range("E14").Value = '  -3.84%  '
$ws.Range("D15").Value = "'" + '0.531'
$ws.Range("E15").Value = '  -4.46%  '
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '27.673.45'
$ws.Range("E16").Value = '  -0.86%  '
$ws.Range("B17").Value = 'Litecoin'
$ws.Range("C17").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D17").Value = "'" + '63.45'
$ws.Range("E17").Value = '  -2.11%  '
$ws.Range("D18").Value = "'" + '219.71'
$ws.Range("E18").Value = '  -3.78%  '
$ws.Range("E19").Value = '  -3.07%  '
$ws.Range("E20").Value = '  -3.64%  '
$ws.Range("E22").Value = '  -4.67%  '
$ws.Range("D23").Value = "'" + '9.68'
$ws.Range("E23").Value = '  -3.28%  '
$ws.Range("E24").Value = '  -3.85%  '
$ws.Range("D25").Value = "'" + '153.98'
$ws.Range("E25").Value = '  -0.57%  '
$ws.Range("D26").Value = "'" + '6.81'
$ws.Range("E26").Value = '  -1.68%  '
$ws.Range("E27").Value = '  +0.19%  '
$ws.Range("E28").Value = '  -1.91%  '
$ws.Range("E29").Value = '  -4.80%  '
$ws.Range("E31").Value = '  -2.28%  '
$ws.Range("E32").Value = '  -5.16%  '
$ws.Range("D33").Value = '1.376.76'
$ws.Range("E33").Value = '  -2.73%  '
$ws.Range("E34").Value = '  -5.12%  '
$ws.Range("E35").Value = '  -4.49%  '
$ws.Range("D36").Value = "'" + '0.975'
$ws.Range("E36").Value = '  -2.38%  '
$ws.Range("E37").Value = '  -0.34%  '
$ws.Range("D38").Value = "'" + '0.0167'
$ws.Range("E38").Value = '  -1.39%  '
$ws.Range("E39").Value = '  -2.94%  '
$ws.Range("D40").Value = "'" + '0.828'
$ws.Range("E40").Value = '  -2.76%  '
$ws.Range("E41").Value = '  +0.19%  '
$ws.Range("E42").Value = '  -3.66%  '
$ws.Range("D43").Value = "'" + '64.39'
$ws.Range("E43").Value = '  -2.09%  '
$ws.Range("E44").Value = '  +2.36%  '
$ws.Range("D45").Value = "'" + '5.21'
$ws.Range("E45").Value = '  -3.64%  '
$ws.Range("E46").Value = '  -5.33%  '
$ws.Range("D47").Value = '1.727.12'
$ws.Range("E47").Value = '  -2.43%  '
$ws.Range("D48").Value = "'" + '87.22'
$ws.Range("E48").Value = '  -1.56%  '
$ws.Range("E49").Value = '  -1.25%  '
$ws.Range("D50").Value = "'" + '0.0967'
$ws.Range("E50").Value = '  -4.07%  '
$ws.Range("E51").Value = '  -1.53%  '
